$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update task estimates in column C
$ws.Range("C3").Value = 3
$ws.Range("C4").Value = 0.5
$ws.Range("C6").Value = 1
$ws.Range("C7").Value = 1.5

# Move the active cell selection from C16 to C12
$ws.Range("C12").Select()

# Force a full recalculation so the SUM formula in C27 (and anything that
# depends on it, e.g. the burndown chart's cached series) is up to date.
$excel.CalculateFullRebuild()
